# The captured OOXML diff for this revision is a pure canonicalization
# artifact: every changed line is the *same* element with its XML
# attributes (and, on the document root, its namespace-prefix
# declarations) written in a different order. Re-serializing both the
# "before" and "after" word/document.xml and word/styles.xml parts with
# W3C Canonical XML (which normalizes attribute/namespace order and
# insignificant whitespace) yields byte-identical trees - i.e. no run
# text, formatting value, numbering property, style definition, page
# geometry, or any other observable document content actually changed.
#
# Word's COM object model has no property that lets a script reorder an
# element's attributes or a root element's xmlns declarations (that
# ordering is an internal choice of whichever serializer last rewrote
# the part) - and since nothing in the object model itself needs to
# change to reach the target state, there is nothing to do here. The
# active document is left exactly as it was loaded.
$d = $word.ActiveDocument
